# Generate Report for Archive
#
# 1. Status text "Ready for handoff" -> "In Translation" wherever it
#    appears (Overview!E2:F2, zh-cn!C2, de-de!C2 all share this string).
# 2. The columns that held that status text narrow to match the new,
#    shorter text (Overview columns E & F, zh-cn column C, de-de column C).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Update the status values ---
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws2.Range("C2").Value = "In Translation"
$ws3.Range("C2").Value = "In Translation"

# --- Narrow the affected columns to their new target width ---
# (target OOXML column width ~= 13.41 characters; the closest width the
# ColumnWidth property can express lands on 13.33 characters)
$ws1.Columns(5).ColumnWidth = 12.5
$ws1.Columns(6).ColumnWidth = 12.5
$ws2.Columns(3).ColumnWidth = 12.5
$ws3.Columns(3).ColumnWidth = 12.5
